{"js": "// Center every table on the page: this sets the table's own horizontal\n// alignment (serialized as <w:jc w:val=\"center\"/> on <w:tblPr>) and, since\n// Word also stamps the chosen alignment onto each row, the alignment of\n// every row in the table too (serialized as <w:jc w:val=\"center\"/> on each\n// <w:trPr>).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const table of tables.items) {\n  table.alignment = Word.Alignment.centered;\n\n  const rows = table.rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  for (const row of rows.items) {\n    row.horizontalAlignment = Word.Alignment.centered;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Center the (single) table on the page: sets table-level alignment\n# (<w:jc w:val=\"center\"/> on <w:tblPr>) AND row-level alignment\n# (<w:jc w:val=\"center\"/> on every <w:trPr>), matching what Word writes\n# when a table's horizontal alignment is set to \"Center\".\n$d = $word.ActiveDocument\n\nforeach ($t in $d.Tables) {\n    # wdAlignRowCenter = 1\n    $t.Alignment = 1\n    foreach ($r in $t.Rows) {\n        $r.Alignment = 1\n    }\n}\n"}
